# "end of first pass through analyses prior to SAF convention"
#
# Extends the residual-basal-area table with a second (currently
# work-in-progress) numeric column in F, mirroring the existing D column
# ("basal_area (sq. ft / acre)") formatting for the rp/control/interior/
# perimeter rows, and leaves the selection parked on the new range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - copy just the number format / alignment (style) of the
# donor cell in D, without disturbing its value, so the new F cells pick up
# exactly the same style index (numFmtId 0.0, centered) that D already uses.
$xlPasteFormats = -4122

# gl / interior  (row 9) - new F cell stays blank, formatted only.
$ws.Range("D9").Copy()
$ws.Range("F9").PasteSpecial($xlPasteFormats)

# gl / perimeter (row 10) - new F cell stays blank, formatted only.
$ws.Range("D10").Copy()
$ws.Range("F10").PasteSpecial($xlPasteFormats)

# rp / interior  (row 12)
$ws.Range("D12").Copy()
$ws.Range("F12").PasteSpecial($xlPasteFormats)
$ws.Range("F12").Value = 8

# rp / perimeter (row 13)
$ws.Range("D13").Copy()
$ws.Range("F13").PasteSpecial($xlPasteFormats)
$ws.Range("F13").Value = 1.5384615384615401

# wedge / interior (row 14) - carries over rp/control's value (18.27586...)
$ws.Range("D14").Copy()
$ws.Range("F14").PasteSpecial($xlPasteFormats)
$ws.Range("F14").Value = 18.275862068965498

$excel.CutCopyMode = $false

# Leave the freshly-populated F12:F14 block selected.
$ws.Range("F12:F14").Select()
